$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1 & 2: "Book Titles" intro paragraph - rewording
# ---------------------------------------------------------------
$d.Content.Find.Execute("must be specified as a number of ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "may be specified as a number of ", 2) | Out-Null

$d.Content.Find.Execute("value pairs. The property", $false, $false, $false, $false, $false, `
    $true, 1, $false, "value pairs. This is necessary where the book titles are not present in the OSIS file. The property", 2) | Out-Null

Write-Output "edit1-2 done"

# ---------------------------------------------------------------
# Edit 3: Replace the "There should be a line..." and "Some OSIS
# Bibles..." paragraphs with five new paragraphs describing the
# BookTitlesInOSIS option.
# ---------------------------------------------------------------
$p1 = $d.Paragraphs(44)
$p2 = $d.Paragraphs(45)
$combined = $d.Range($p1.Range.Start, $p2.Range.End)

$para1 = "Some OSIS Bibles include book titles as headings in the book introduction or as initial section headings.  Where this is the case, book titles may be extracted in from the OSIS file by including the line:"
$para2 = "BookTitlesInOSIS=True"
$para3 = "Where this is specified, a title at that start of a book is assumed to be the book title. If titles are included in the OSIS file only for some books, titles for the other books may still be specified in convert.txt, but book titles extracted from the OSIS file take priority."
$para4 = "There are some cases where it may be necessary to specify book titles in convert.txt even though they are present in the OSIS file. This will be the case, for example, where a Bible or Testament introduction precedes the book title for the first book of the Bible or Testament.  The titles specified in convert.txt should match exactly the titles in the OSIS file; the conversion tool will detect that the titles match and so will recognise the book title even where it does not occur at the start of the book."
$para5 = "If there is any case where a book starts with a title which in not the book title, “BookTitlesInOSIS=True” should not be specified, as in the case where the book title is preceded by a Bible or Testament introduction which has a title. If titles in convert.txt  match exactly the book titles in the OSIS file, the conversion tool ensure that duplicate titles do not appear."

$combined.Text = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5

Write-Output "edit3 done"
